$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-26 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-27 Monday", 1) | Out-Null
$d.Content.Find.Execute("33-19=", $true, $false, $false, $false, $false, $true, 1, $false, "70+25=", 1) | Out-Null
$d.Content.Find.Execute("2+13=", $true, $false, $false, $false, $false, $true, 1, $false, "2+86=", 1) | Out-Null
$d.Content.Find.Execute("86-15=", $true, $false, $false, $false, $false, $true, 1, $false, "52+32=", 1) | Out-Null
$d.Content.Find.Execute("33+30=", $true, $false, $false, $false, $false, $true, 1, $false, "90-53=", 1) | Out-Null
$d.Content.Find.Execute("8+72=", $true, $false, $false, $false, $false, $true, 1, $false, "49+34=", 1) | Out-Null
$d.Content.Find.Execute("42-33=", $true, $false, $false, $false, $false, $true, 1, $false, "21-3=", 1) | Out-Null
$d.Content.Find.Execute("15-0=", $true, $false, $false, $false, $false, $true, 1, $false, "13+31=", 1) | Out-Null
$d.Content.Find.Execute("13+57=", $true, $false, $false, $false, $false, $true, 1, $false, "49+11=", 1) | Out-Null
$d.Content.Find.Execute("11+21=", $true, $false, $false, $false, $false, $true, 1, $false, "8+73=", 1) | Out-Null
$d.Content.Find.Execute("85+5=", $true, $false, $false, $false, $false, $true, 1, $false, "82-56=", 1) | Out-Null
$d.Content.Find.Execute("92-54=", $true, $false, $false, $false, $false, $true, 1, $false, "78+15=", 1) | Out-Null
$d.Content.Find.Execute("32-17=", $true, $false, $false, $false, $false, $true, 1, $false, "45-35=", 1) | Out-Null
$d.Content.Find.Execute("19+67=", $true, $false, $false, $false, $false, $true, 1, $false, "97-36=", 1) | Out-Null
$d.Content.Find.Execute("0+42=", $true, $false, $false, $false, $false, $true, 1, $false, "63-11=", 1) | Out-Null
$d.Content.Find.Execute("12+69=", $true, $false, $false, $false, $false, $true, 1, $false, "57-8=", 1) | Out-Null
$d.Content.Find.Execute("31+37=", $true, $false, $false, $false, $false, $true, 1, $false, "81+1=", 1) | Out-Null
$d.Content.Find.Execute("52+47=", $true, $false, $false, $false, $false, $true, 1, $false, "82-43=", 1) | Out-Null
$d.Content.Find.Execute("60-42=", $true, $false, $false, $false, $false, $true, 1, $false, "51-37=", 1) | Out-Null
$d.Content.Find.Execute("76-72=", $true, $false, $false, $false, $false, $true, 1, $false, "24+37=", 1) | Out-Null
$d.Content.Find.Execute("90-57=", $true, $false, $false, $false, $false, $true, 1, $false, "55+38=", 1) | Out-Null
$d.Content.Find.Execute("80-48=", $true, $false, $false, $false, $false, $true, 1, $false, "50+0=", 1) | Out-Null
$d.Content.Find.Execute("61+14=", $true, $false, $false, $false, $false, $true, 1, $false, "2-1=", 1) | Out-Null
$d.Content.Find.Execute("37-5=", $true, $false, $false, $false, $false, $true, 1, $false, "87+7=", 1) | Out-Null
$d.Content.Find.Execute("89-29=", $true, $false, $false, $false, $false, $true, 1, $false, "34+28=", 1) | Out-Null
$d.Content.Find.Execute("8+38=", $true, $false, $false, $false, $false, $true, 1, $false, "24-18=", 1) | Out-Null
$d.Content.Find.Execute("61-56=", $true, $false, $false, $false, $false, $true, 1, $false, "61-17=", 1) | Out-Null
$d.Content.Find.Execute("35+30=", $true, $false, $false, $false, $false, $true, 1, $false, "75-10=", 1) | Out-Null
$d.Content.Find.Execute("24-2=", $true, $false, $false, $false, $false, $true, 1, $false, "42-3=", 1) | Out-Null
$d.Content.Find.Execute("52-25=", $true, $false, $false, $false, $false, $true, 1, $false, "85+0=", 1) | Out-Null
$d.Content.Find.Execute("98-87=", $true, $false, $false, $false, $false, $true, 1, $false, "40+42=", 1) | Out-Null
$d.Content.Find.Execute("82+2=", $true, $false, $false, $false, $false, $true, 1, $false, "87-1=", 1) | Out-Null
$d.Content.Find.Execute("28+4=", $true, $false, $false, $false, $false, $true, 1, $false, "92+1=", 1) | Out-Null
$d.Content.Find.Execute("25+27=", $true, $false, $false, $false, $false, $true, 1, $false, "41+43=", 1) | Out-Null
$d.Content.Find.Execute("75+0=", $true, $false, $false, $false, $false, $true, 1, $false, "9+24=", 1) | Out-Null
$d.Content.Find.Execute("87+0=", $true, $false, $false, $false, $false, $true, 1, $false, "16+82=", 1) | Out-Null
$d.Content.Find.Execute("87-48=", $true, $false, $false, $false, $false, $true, 1, $false, "15+33=", 1) | Out-Null
$d.Content.Find.Execute("54+22=", $true, $false, $false, $false, $false, $true, 1, $false, "34-12=", 1) | Out-Null
$d.Content.Find.Execute("76-74=", $true, $false, $false, $false, $false, $true, 1, $false, "66+5=", 1) | Out-Null
$d.Content.Find.Execute("84+1=", $true, $false, $false, $false, $false, $true, 1, $false, "2+73=", 1) | Out-Null
$d.Content.Find.Execute("75-20=", $true, $false, $false, $false, $false, $true, 1, $false, "10+54=", 1) | Out-Null
$d.Content.Find.Execute("96-86=", $true, $false, $false, $false, $false, $true, 1, $false, "96-8=", 1) | Out-Null
$d.Content.Find.Execute("37+21=", $true, $false, $false, $false, $false, $true, 1, $false, "26+50=", 1) | Out-Null
$d.Content.Find.Execute("90-30=", $true, $false, $false, $false, $false, $true, 1, $false, "60-7=", 1) | Out-Null
$d.Content.Find.Execute("73-31=", $true, $false, $false, $false, $false, $true, 1, $false, "39-22=", 1) | Out-Null
$d.Content.Find.Execute("50+11=", $true, $false, $false, $false, $false, $true, 1, $false, "50-49=", 1) | Out-Null
$d.Content.Find.Execute("34-1=", $true, $false, $false, $false, $false, $true, 1, $false, "60-18=", 1) | Out-Null
$d.Content.Find.Execute("83-32=", $true, $false, $false, $false, $false, $true, 1, $false, "76-47=", 1) | Out-Null
$d.Content.Find.Execute("24+74=", $true, $false, $false, $false, $false, $true, 1, $false, "14+8=", 1) | Out-Null
$d.Content.Find.Execute("20+30=", $true, $false, $false, $false, $false, $true, 1, $false, "67-20=", 1) | Out-Null
$d.Content.Find.Execute("84-67=", $true, $false, $false, $false, $false, $true, 1, $false, "83-50=", 1) | Out-Null
$d.Content.Find.Execute("72-31=", $true, $false, $false, $false, $false, $true, 1, $false, "35+9=", 1) | Out-Null
$d.Content.Find.Execute("40+7=", $true, $false, $false, $false, $false, $true, 1, $false, "75-58=", 1) | Out-Null
$d.Content.Find.Execute("66-13=", $true, $false, $false, $false, $false, $true, 1, $false, "98-96=", 1) | Out-Null
$d.Content.Find.Execute("89-28=", $true, $false, $false, $false, $false, $true, 1, $false, "56-42=", 1) | Out-Null
$d.Content.Find.Execute("85-74=", $true, $false, $false, $false, $false, $true, 1, $false, "14+75=", 1) | Out-Null
$d.Content.Find.Execute("82-27=", $true, $false, $false, $false, $false, $true, 1, $false, "28+55=", 1) | Out-Null
$d.Content.Find.Execute("38+17=", $true, $false, $false, $false, $false, $true, 1, $false, "11+34=", 1) | Out-Null
$d.Content.Find.Execute("56-51=", $true, $false, $false, $false, $false, $true, 1, $false, "54+33=", 1) | Out-Null
$d.Content.Find.Execute("48-32=", $true, $false, $false, $false, $false, $true, 1, $false, "69-41=", 1) | Out-Null
$d.Content.Find.Execute("47+35=", $true, $false, $false, $false, $false, $true, 1, $false, "9+65=", 1) | Out-Null
$d.Content.Find.Execute("82-2=", $true, $false, $false, $false, $false, $true, 1, $false, "51-0=", 1) | Out-Null
$d.Content.Find.Execute("27+4=", $true, $false, $false, $false, $false, $true, 1, $false, "49-40=", 1) | Out-Null
$d.Content.Find.Execute("58-5=", $true, $false, $false, $false, $false, $true, 1, $false, "42+1=", 1) | Out-Null
$d.Content.Find.Execute("30+28=", $true, $false, $false, $false, $false, $true, 1, $false, "15-2=", 1) | Out-Null
$d.Content.Find.Execute("79-55=", $true, $false, $false, $false, $false, $true, 1, $false, "5+22=", 1) | Out-Null
$d.Content.Find.Execute("13+2=", $true, $false, $false, $false, $false, $true, 1, $false, "99-24=", 1) | Out-Null
$d.Content.Find.Execute("74-55=", $true, $false, $false, $false, $false, $true, 1, $false, "64-21=", 1) | Out-Null
$d.Content.Find.Execute("52-15=", $true, $false, $false, $false, $false, $true, 1, $false, "51-43=", 1) | Out-Null
$d.Content.Find.Execute("93-29=", $true, $false, $false, $false, $false, $true, 1, $false, "84-32=", 1) | Out-Null
$d.Content.Find.Execute("5+94=", $true, $false, $false, $false, $false, $true, 1, $false, "36+55=", 1) | Out-Null
$d.Content.Find.Execute("87-75=", $true, $false, $false, $false, $false, $true, 1, $false, "80-28=", 1) | Out-Null
$d.Content.Find.Execute("21+77=", $true, $false, $false, $false, $false, $true, 1, $false, "99-12=", 1) | Out-Null
$d.Content.Find.Execute("79-52=", $true, $false, $false, $false, $false, $true, 1, $false, "86-46=", 1) | Out-Null
$d.Content.Find.Execute("29-10=", $true, $false, $false, $false, $false, $true, 1, $false, "24-23=", 1) | Out-Null
$d.Content.Find.Execute("26+23=", $true, $false, $false, $false, $false, $true, 1, $false, "37+6=", 1) | Out-Null
$d.Content.Find.Execute("43-39=", $true, $false, $false, $false, $false, $true, 1, $false, "3+27=", 1) | Out-Null
$d.Content.Find.Execute("84-53=", $true, $false, $false, $false, $false, $true, 1, $false, "47+22=", 1) | Out-Null
$d.Content.Find.Execute("49-45=", $true, $false, $false, $false, $false, $true, 1, $false, "44+20=", 1) | Out-Null
$d.Content.Find.Execute("32-13=", $true, $false, $false, $false, $false, $true, 1, $false, "14+39=", 1) | Out-Null
$d.Content.Find.Execute("1+91=", $true, $false, $false, $false, $false, $true, 1, $false, "41-23=", 1) | Out-Null
$d.Content.Find.Execute("4+6=", $true, $false, $false, $false, $false, $true, 1, $false, "12+17=", 1) | Out-Null
$d.Content.Find.Execute("6+45=", $true, $false, $false, $false, $false, $true, 1, $false, "69-14=", 1) | Out-Null
$d.Content.Find.Execute("60-10=", $true, $false, $false, $false, $false, $true, 1, $false, "8+77=", 1) | Out-Null
$d.Content.Find.Execute("93-58=", $true, $false, $false, $false, $false, $true, 1, $false, "83-22=", 1) | Out-Null
$d.Content.Find.Execute("33-15=", $true, $false, $false, $false, $false, $true, 1, $false, "61-55=", 1) | Out-Null
$d.Content.Find.Execute("80-3=", $true, $false, $false, $false, $false, $true, 1, $false, "83-49=", 1) | Out-Null
$d.Content.Find.Execute("60-42=", $true, $false, $false, $false, $false, $true, 1, $false, "14+63=", 1) | Out-Null
$d.Content.Find.Execute("23-5=", $true, $false, $false, $false, $false, $true, 1, $false, "75-2=", 1) | Out-Null
$d.Content.Find.Execute("36-14=", $true, $false, $false, $false, $false, $true, 1, $false, "57-15=", 1) | Out-Null
$d.Content.Find.Execute("48+47=", $true, $false, $false, $false, $false, $true, 1, $false, "11+63=", 1) | Out-Null
$d.Content.Find.Execute("67-34=", $true, $false, $false, $false, $false, $true, 1, $false, "29+57=", 1) | Out-Null
$d.Content.Find.Execute("26+62=", $true, $false, $false, $false, $false, $true, 1, $false, "0+20=", 1) | Out-Null
$d.Content.Find.Execute("84-82=", $true, $false, $false, $false, $false, $true, 1, $false, "15+23=", 1) | Out-Null
$d.Content.Find.Execute("30-11=", $true, $false, $false, $false, $false, $true, 1, $false, "84-11=", 1) | Out-Null
$d.Content.Find.Execute("25+21=", $true, $false, $false, $false, $false, $true, 1, $false, "88-36=", 1) | Out-Null
$d.Content.Find.Execute("69+1=", $true, $false, $false, $false, $false, $true, 1, $false, "41+36=", 1) | Out-Null
$d.Content.Find.Execute("98-56=", $true, $false, $false, $false, $false, $true, 1, $false, "84-54=", 1) | Out-Null
$d.Content.Find.Execute("56+8=", $true, $false, $false, $false, $false, $true, 1, $false, "20+21=", 1) | Out-Null
$d.Content.Find.Execute("36+20=", $true, $false, $false, $false, $false, $true, 1, $false, "92-3=", 1) | Out-Null
$d.Content.Find.Execute("34+5=", $true, $false, $false, $false, $false, $true, 1, $false, "77-13=", 1) | Out-Null
